# Weekly update: shift existing "Alcachofa" price rows down to make room for new
# observations at the top of the time series, and append newest weekly rows at the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{}
$data[65] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45126, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 40, 17000, 17000, 17000, "`$/caja 50 unidades", "Provincia de Limarí", 340, 50, "Hortaliza")
$data[66] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45126, 16, 100112013, "Alcachofa", "Española", "Primera", 50, 17000, 17000, 17000, "`$/caja 30 unidades", "Provincia de Limarí", 567, 30, "Hortaliza")
$data[67] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45113, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 60, 16000, 16000, 16000, "`$/caja 50 unidades", "Provincia de Limarí", 320, 50, "Hortaliza")
$data[68] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44775, 16, 100112013, "Alcachofa", "Española", "Primera", 60, 18000, 19000, 18500, "`$/caja 30 unidades", "Provincia de Limarí", 617, 30, "Hortaliza")
$data[69] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44775, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 14000, 15000, 14500, "`$/caja 40 unidades", "Provincia de Limarí", 362, 40, "Hortaliza")
$data[70] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44755, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 60, 16000, 17000, 16500, "`$/caja 50 unidades", "Provincia de Limarí", 330, 50, "Hortaliza")
$data[71] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44508, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 160, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia del Elquí", 288, 40, "Hortaliza")
$data[72] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44816, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 60, 13000, 14000, 13500, "`$/caja 50 unidades", "Provincia de Limarí", 270, 50, "Hortaliza")
$data[73] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44789, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 100, 14000, 15000, 14500, "`$/caja 50 unidades", "Provincia de Limarí", 290, 50, "Hortaliza")
$data[74] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44789, 16, 100112013, "Alcachofa", "Española", "Primera", 80, 15000, 16000, 15500, "`$/caja 30 unidades", "Provincia de Limarí", 517, 30, "Hortaliza")
$data[75] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44789, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 13000, 14000, 13500, "`$/caja 40 unidades", "Provincia de Limarí", 338, 40, "Hortaliza")
$data[76] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44813, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 80, 13000, 14000, 13500, "`$/caja 50 unidades", "Provincia de Limarí", 270, 50, "Hortaliza")
$data[77] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44813, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 100, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia de Limarí", 288, 40, "Hortaliza")
$data[78] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44875, 16, 100112013, "Alcachofa", "Española", "Primera", 60, 10000, 10000, 10000, "`$/caja 30 unidades", "Provincia de Limarí", 333, 30, "Hortaliza")
$data[79] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44516, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 120, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia del Elquí", 288, 40, "Hortaliza")
$data[80] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44435, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 120, 14000, 15000, 14500, "`$/caja 40 unidades", "Provincia del Elquí", 362, 40, "Hortaliza")
$data[81] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45120, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 17000, 17000, 17000, "`$/caja 40 unidades", "Provincia de Limarí", 425, 40, "Hortaliza")
$data[82] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44831, 16, 100112013, "Alcachofa", "Española", "Primera", 60, 11000, 12000, 11500, "`$/caja 30 unidades", "Provincia de Limarí", 383, 30, "Hortaliza")
$data[83] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44831, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 10000, 10000, 10000, "`$/caja 40 unidades", "Provincia de Limarí", 250, 40, "Hortaliza")
$data[84] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44473, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 160, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia del Elquí", 288, 40, "Hortaliza")
$data[85] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44427, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 120, 13000, 14000, 13500, "`$/caja 40 unidades", "Provincia del Elquí", 338, 40, "Hortaliza")
$data[86] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44503, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 160, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia del Elquí", 288, 40, "Hortaliza")
$data[87] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45097, 16, 100112013, "Alcachofa", "Española", "Primera", 120, 15000, 16000, 15500, "`$/caja 30 unidades", "Provincia de Limarí", 517, 30, "Hortaliza")
$data[88] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44784, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 60, 14000, 15000, 14500, "`$/caja 50 unidades", "Provincia de Limarí", 290, 50, "Hortaliza")
$data[89] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44784, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 13000, 14000, 13500, "`$/caja 40 unidades", "Provincia de Limarí", 338, 40, "Hortaliza")
$data[90] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45085, 16, 100112013, "Alcachofa", "Española", "Primera", 20, 16000, 16000, 16000, "`$/caja 50 unidades", "Provincia de Limarí", 320, 50, "Hortaliza")
$data[91] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44838, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 100, 9000, 10000, 9500, "`$/caja 50 unidades", "Provincia de Limarí", 190, 50, "Hortaliza")
$data[92] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44838, 16, 100112013, "Alcachofa", "Española", "Primera", 100, 10000, 11000, 10500, "`$/caja 30 unidades", "Provincia de Limarí", 350, 30, "Hortaliza")
$data[93] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44838, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 100, 9000, 10000, 9500, "`$/caja 40 unidades", "Provincia de Limarí", 238, 40, "Hortaliza")
$data[94] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45079, 16, 100112013, "Alcachofa", "Española", "Primera", 30, 16000, 16000, 16000, "`$/caja 30 unidades", "Provincia de Limarí", 533, 30, "Hortaliza")
$data[95] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44490, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 100, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia del Elquí", 288, 40, "Hortaliza")
$data[96] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45124, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 50, 17000, 17000, 17000, "`$/caja 50 unidades", "Provincia de Limarí", 340, 50, "Hortaliza")
$data[97] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44495, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 120, 11000, 12000, 11500, "`$/caja 40 unidades", "Provincia del Elquí", 288, 40, "Hortaliza")
$data[98] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44417, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 120, 15000, 16000, 15500, "`$/caja 40 unidades", "Provincia del Elquí", 388, 40, "Hortaliza")
$data[99] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44468, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 12000, 13000, 12500, "`$/caja 40 unidades", "Provincia del Elquí", 312, 40, "Hortaliza")
$data[100] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44420, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 120, 13000, 14000, 13500, "`$/caja 40 unidades", "Provincia del Elquí", 338, 40, "Hortaliza")
$data[101] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44811, 16, 100112013, "Alcachofa", "Española", "Primera", 60, 12000, 13000, 12500, "`$/caja 30 unidades", "Provincia de Limarí", 417, 30, "Hortaliza")
$data[102] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45112, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 50, 16000, 16000, 16000, "`$/caja 50 unidades", "Provincia de Limarí", 320, 50, "Hortaliza")
$data[103] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45089, 16, 100112013, "Alcachofa", "Argentina(o)", "Primera", 60, 16000, 17000, 16500, "`$/caja 50 unidades", "Provincia de Limarí", 330, 50, "Hortaliza")
$data[104] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 45121, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 40, 17000, 17000, 17000, "`$/caja 40 unidades", "Provincia de Limarí", 425, 40, "Hortaliza")
$data[105] = @(7, "Terminal Hortofrutícola Agro Chillán", "Ñuble", 44777, 16, 100112013, "Alcachofa", "Madrigal", "Primera", 60, 14000, 15000, 14500, "`$/caja 40 unidades", "Provincia del Elquí", 362, 40, "Hortaliza")

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
    # Column D holds the observation date; keep the sheet-wide date format.
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Write-Host "Updated rows 65-105 (dimension now A1:R105)"
